# "a bunch of runs" - adds percent-diff comparison columns (D:G) for the
# first ensemble table, converts several duplicate single-cell formulas
# into shared formula ranges, and refreshes some cosmetic formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Restyle the "Number" column (B3:B12, B19:B28) - keep centered.
#    (value/meaning unchanged - Excel simply re-centers on resave)
# ---------------------------------------------------------------------
$ws.Range("B3:B12").HorizontalAlignment = -4108
$ws.Range("B19:B28").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# 2) Restyle "Initial Drawdown Aq" column (C3:C12): black Calibri,
#    centered + vertically centered.
# ---------------------------------------------------------------------
foreach ($addr in @("C3", "C4", "C5", "C6", "C10", "C12")) {
    $r = $ws.Range($addr)
    $r.Font.Color = 0
    $r.HorizontalAlignment = -4108
    $r.VerticalAlignment = -4108
}
foreach ($addr in @("C7", "C8", "C9", "C11")) {
    $r = $ws.Range($addr)
    $r.HorizontalAlignment = -4108
    $r.VerticalAlignment = -4108
}

# ---------------------------------------------------------------------
# 3) New comparison columns D:G for rows 3-12
#    D = "Drawdown crop change", E = "% Diff", F = "Drawdown Well Move",
#    G = "%Diff" (header row already has these labels in row 2)
# ---------------------------------------------------------------------

# D3:D4 / E3:E4 stay empty but take on the centered format
$ws.Range("D3:E4").HorizontalAlignment = -4108
$ws.Range("D3:E4").VerticalAlignment = -4108

$dVals = @{
    5  = -0.53250120000000001
    6  = -9.6519659999999998
    7  = -2.0276565999999998
    8  = -2.0350646999999999
    9  = -2.0354766999999998
    10 = -0.96211239999999998
    11 = -1.2331467
    12 = -9.7687840000000001
}
foreach ($row in $dVals.Keys) {
    $ws.Range("D$row").Value = $dVals[$row]
}

$ws.Range("E5:E11").Formula = "=(C5-D5) /C5 * 100"
$ws.Range("E5:E11").HorizontalAlignment = -4108
$ws.Range("E5:E11").VerticalAlignment = -4108

$ws.Range("E12").Formula = "=(C12-D12) /C12 * 100"
$ws.Range("E12").HorizontalAlignment = -4108
$ws.Range("E12").VerticalAlignment = -4108

$fVals = @{
    3  = -4.3836975000000002
    4  = -0.48382567999999998
    5  = -0.56997679999999995
    6  = -9.765053
    7  = -2.1143339999999999
    8  = -2.1297760000000001
    9  = -2.1286315999999998
    10 = -0.95437620000000001
    11 = -1.3412476
    12 = -9.9868509999999997
}
foreach ($row in $fVals.Keys) {
    $ws.Range("F$row").Value = $fVals[$row]
}
$ws.Range("F3:F4").Font.Color = 0
$ws.Range("F3:F4").HorizontalAlignment = -4108
$ws.Range("F3:F4").VerticalAlignment = -4108
$ws.Range("F5:F12").HorizontalAlignment = -4108
$ws.Range("F5:F12").VerticalAlignment = -4108

$ws.Range("G3").Formula = "=ABS(((F3-C3)/F3)*100)"
$ws.Range("G4:G12").Formula = "=ABS(((F4-C4)/F4)*100)"
$ws.Range("G3:G12").HorizontalAlignment = -4108
$ws.Range("G3:G12").VerticalAlignment = -4108

# ---------------------------------------------------------------------
# 4) Update the second table's D12 input + consolidate its per-row
#    "% change" formulas (E22:E28) into one shared formula.
# ---------------------------------------------------------------------
$ws.Range("E22:E28").Formula = "=(1-D22/C22)*100"

# Column I (old scratch formatting) is no longer needed for rows 26-28.
$ws.Range("I26:I28").Clear()
$ws.Range("I22:I25").Font.Name = "Courier New"
$ws.Range("I22:I25").Font.Size = 14
$ws.Range("I22:I25").Font.Color = 0

# ---------------------------------------------------------------------
# 5) Decorative empty-cell formatting left over from selecting ranges
#    K7:K16 and M2:N12 while picking fonts (no values were entered).
# ---------------------------------------------------------------------
$ws.Range("K7:K15").Font.Name = "Courier New"
$ws.Range("K7:K15").Font.Size = 10
$ws.Range("K7:K15").Font.Color = 0

$ws.Range("K16").Font.Name = "Courier New"
$ws.Range("K16").Font.Size = 14
$ws.Range("K16").Font.Color = 3159199

$ws.Range("M2:M12").Font.Name = "Courier New"
$ws.Range("M2:M12").Font.Size = 14
$ws.Range("M2:M12").Font.Color = 0

$ws.Range("N2").Font.Name = "Courier New"
$ws.Range("N2").Font.Size = 14
$ws.Range("N2").Font.Color = 0

# ---------------------------------------------------------------------
# 6) Row-height refresh (matches the resave with a newer Excel build)
# ---------------------------------------------------------------------
$ws.Range("A1").RowHeight = 16
$ws.Range("A2").RowHeight = 45
$ws.Range("A3:A12").RowHeight = 20
$ws.Range("A16").RowHeight = 19
$ws.Range("A17:A21").RowHeight = 16
$ws.Range("A22:A25").RowHeight = 20
$ws.Range("A26:A28").RowHeight = 16

# ---------------------------------------------------------------------
# 7) Selection, matching the author's last on-screen state.
# ---------------------------------------------------------------------
$ws.Range("J7:N16").Select()

$excel.Calculate()
